$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 416: updated new-positive-cases count (B416's cumulative-sum formula
# recalculates automatically from this, cascading into the following rows).
$ws.Range("C416").Value = 122

# Row 418: updated new-positive-cases count.
$ws.Range("C418").Value = 37

# Row 419: updated new-positive-cases count.
$ws.Range("C419").Value = 122

# Row 420: this day's row was still an empty placeholder (formulas showing
# blank "" results) - fill in the real daily figures that were published.
$ws.Range("C420").Value = 11
$ws.Range("E420").Value = 8
$ws.Range("F420").Value = 8
$ws.Range("G420").Value = 37

# L420/M420 are formatted as Text ("@"), so writing straight into .Value
# would store the literal string "0" rather than the number 0 (matching
# how Excel itself treats typed input in a Text-formatted cell). Swap to
# General just long enough to write the real number, then restore the
# Text format so the cell's formatting stays exactly as it was.
$ws.Range("L420").NumberFormat = "General"
$ws.Range("L420").Value = 0
$ws.Range("L420").NumberFormat = "@"

$ws.Range("M420").NumberFormat = "General"
$ws.Range("M420").Value = 0
$ws.Range("M420").NumberFormat = "@"
